# Database_tables.xlsx edit script
# "everything ready except Testing"
#
# Applies the content changes described by the commit diff:
#  - "pass" column: type varchar(50) -> varchar(100); comment "hashed?" ->
#    "hashed value using sha256"
#  - "alias" (accountXref) comment wording: "use to allow..." -> "used to
#    allow..."
#  - new comments added for the phone/email columns in the "info" table
#  - the old "*transactions logged, table not necessary" comment is cleared
#  - the running concatenated-SQL helper formula (E.. / F..) for the "info"
#    table block moved down one row (was anchored on the blank header row,
#    now anchored on the first real column row) with a trailing space added
#    to match the other blocks
#  - the u_id/first_name/.../pass "UPDATE" helper formulas for the "users"
#    table gained a trailing space after "?," so the concatenated clause
#    list reads "a = ?, b = ?, " instead of "a = ?,b = ?,"
#  - cosmetic: column D narrower, selection moved to D29, no frozen/ scrolled
#    topLeftCell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- "users" table UPDATE-clause helper column (column E, rows 4-9) -------
# Previously: =CONCAT(B4," = ?,")   Now: =CONCAT(B4," = ?, ")  (trailing space)
$ws.Range("E4").Formula = "=CONCAT(B4, "" = ?, "")"
$ws.Range("E5").Formula = "=CONCAT(B5, "" = ?, "")"
$ws.Range("E6").Formula = "=CONCAT(B6, "" = ?, "")"
$ws.Range("E7").Formula = "=CONCAT(B7, "" = ?, "")"
$ws.Range("E8").Formula = "=CONCAT(B8, "" = ?, "")"
$ws.Range("E9").Formula = "=CONCAT(B9, "" = ?, "")"

# --- "pass" row: widen the stored type and clarify the comment -----------
$ws.Range("C9").Value = "varchar(100)"
$ws.Range("D9").Value = "hashed value using sha256"

# --- "information" table block: drop the old header-row concat cells and -
#     recreate the full-clause concat one row down, on the first real field
$ws.Range("E20").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Range("F21").Formula = "=CONCAT(E21:E29)"

# --- new explanatory comments for phone / email ---------------------------
$ws.Range("D27").Value = "special characters removed before insertion and then added after retrieval"
$ws.Range("D28").Value = "regular expressions used to verify format before insertion"

# --- accountXref comment: reword "use to allow" -> "used to allow" --------
$ws.Range("D30").Value = "used to allow for users to have multiple accounts & joint accounts"

# --- stale comment removed from the bottom notes area ----------------------
$ws.Range("D35").ClearContents()

# --- cosmetic: narrower comment column, scrolled/selected cell ------------
$ws.Columns.Item(4).ColumnWidth = 80

$ws.Range("D29").Select() | Out-Null
